$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 82
$ws1.Range("F3").Value = 809
$ws1.Range("F6").Value = 113
$ws1.Range("F8").Value = 4510
$ws1.Range("F9").Value = 98
$ws1.Range("F10").Value = 5012
$ws1.Range("F11").Value = 566
$ws1.Range("F12").Value = 1265

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 82
$ws4.Range("F3").Value = 809
$ws4.Range("F6").Value = 113
$ws4.Range("F9").Value = 4510
$ws4.Range("F10").Value = 98
$ws4.Range("F11").Value = 5012
$ws4.Range("F12").Value = 566
$ws4.Range("F13").Value = 1265
